$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are stored as text in the workbook (inline strings). Some of the
# new "Price" values look like plain numbers to Excel (e.g. "1.000", "41.59"),
# so a leading apostrophe is used to force them to stay text, matching the
# original text-typed cells, instead of being auto-converted to numeric values.

$ws.Range("D2").Value = '28.122.49'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '1.874.52'
$ws.Range("E3").Value = '  -1.63%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''313.09'
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").Value = '''0.9996'
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = '''0.5050'
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").Value = '''0.3846'
$ws.Range("E8").Value = '  -2.33%  '
$ws.Range("D9").Value = '''0.09040'
$ws.Range("E9").Value = '  -5.94%  '
$ws.Range("D10").Value = '''1.122'
$ws.Range("E10").Value = '  -1.38%  '
$ws.Range("D11").Value = '''41.59'
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("D12").Value = '''6.364'
$ws.Range("E12").Value = '  -0.50%  '
$ws.Range("D13").Value = '''20.77'
$ws.Range("E13").Value = '  -0.54%  '
$ws.Range("D14").Value = '1.878.58'
$ws.Range("E14").Value = '  -0.88%  '
$ws.Range("D15").Value = '''7.277'
$ws.Range("E15").Value = '  -0.84%  '
$ws.Range("D16").Value = '''1.001'
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").Value = '''0.00001109'
$ws.Range("E17").Value = '  -1.08%  '
$ws.Range("D18").Value = '''91.30'
$ws.Range("E18").Value = '  -1.16%  '
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").Value = '''18.21'
$ws.Range("E20").Value = '  +1.55%  '
$ws.Range("D21").Value = '''0.9993'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("E22").Value = '  -1.15%  '
$ws.Range("D23").Value = '28.150.01'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  +1.44%  '
$ws.Range("D25").Value = '''2.260'
$ws.Range("E25").Value = '  -1.93%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").Value = '''3.399'
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '''2.548'
$ws.Range("E27").Value = '  -4.10%  '
$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").Value = '2.091.85'
$ws.Range("E28").Value = '  -0.84%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '''20.84'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = '''156.92'
$ws.Range("E30").Value = '  -0.44%  '
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = '''126.99'
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '''0.1064'
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''1.063'
$ws.Range("E33").Value = '  -2.32%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '''5.617'
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''3.599'
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = '''9.480'
$ws.Range("E36").Value = '  -0.83%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '''0.06594'
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.02408'
$ws.Range("E38").Value = '  -1.00%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = '''0.2195'
$ws.Range("E39").Value = '  +0.54%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '''1.297'
$ws.Range("E40").Value = '  +2.94%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '''1.213'
$ws.Range("E41").Value = '  -1.93%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '''0.6409'
$ws.Range("E42").Value = '  +0.94%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '''11.51'
$ws.Range("E43").Value = '  +1.35%  '
$ws.Range("B44").Value = 'InternetComputer(DFINITY)'
$ws.Range("C44").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D44").Value = '''4.930'
$ws.Range("E44").Value = '  -1.67%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").Value = '''0.9996'
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''13.26'
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.6034'
$ws.Range("E47").Value = '  +0.59%  '
$ws.Range("B48").Value = 'WEMIXTOKEN'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '''1.275'
$ws.Range("E48").Value = '  -0.22%  '
$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D49").Value = '''3.667'
$ws.Range("E49").Value = '  -1.60%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = '''1.236'
$ws.Range("E50").Value = '  +4.32%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '''2.008'
$ws.Range("E51").Value = '  -0.93%  '
